$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SystemSettings")

# Add new shared string value used by output formatter for templateExport
$ws.Range("C2").Value = '$(Data:Enabled)'

# Update the active selection on the sheet to reflect where the user left off
$ws.Range("C3").Select()
